# präsi + folien aktualisiert
#
# Slide 6 ("Ergebnisse"): the content placeholder text that spelled out the
# homepage URL is replaced by a short "STARTSEITE" label that now carries an
# actual hyperlink to the HydroTirol GitHub Pages site.

$p = $ppt.ActivePresentation

# --- Slide 6: turn the plain text into a hyperlinked "STARTSEITE" label ---
$slide = $p.Slides.Item(6)

# Find the content placeholder that currently holds the spelled-out URL
# (fall back to the known placeholder shape if the text can't be matched).
$shape = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $candidate = $slide.Shapes.Item($i)
    if ($candidate.HasTextFrame -and $candidate.TextFrame.HasText) {
        if ($candidate.TextFrame.TextRange.Text -match "Startseite") {
            $shape = $candidate
            break
        }
    }
}
if ($null -eq $shape) {
    $shape = $slide.Shapes.Item(2)
}

$textRange = $shape.TextFrame.TextRange
$textRange.Text = "STARTSEITE"
$hyperlink = $textRange.ActionSettings(1).Hyperlink
$hyperlink.Address = "https://hydrotirol.github.io/index.html"

# --- Presentation-level slide-guide bookkeeping (empty guide list marker) ---
# Touching the guides collection is what stamps the (empty) p15:sldGuideLst
# extension PowerPoint writes into presentation.xml once the Guides feature
# has been used.
try {
    $guides = $p.Guides
    $newGuide = $guides.Add(1, 3.0)
    $newGuide.Delete()
} catch {
    # Guides automation may be unavailable in some hosts; ignore if so.
}
